$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 41: 23rd April 2020 update ---

# A41: Date (serial 43944 = 2020-04-23), keep the same date formatting/alignment
# used by the other date cells in column A (style index 2: numFmt 15, centered).
$ws.Range("A41").Value = 43944
$ws.Range("A41").NumberFormat = "d-mmm-yy"
$ws.Range("A41").HorizontalAlignment = -4108  # xlCenter

# B41: New Cases
$ws.Range("B41").Value = 17

# C41: Tested
$ws.Range("C41").Value = 668

# D41: Travelled From
$ws.Range("D41").Value = "None"

# E41: County
$ws.Range("E41").Value = "Mombasa(12),Nairobi(3)"

# F41: Aggregation
$ws.Range("F41").Value = 320

# G41: Case Type
$ws.Range("G41").Value = "Community(17)"

# H41: Recover
$ws.Range("H41").Value = 0

# I41: Death
$ws.Range("I41").Value = 0

# K41: Info Giver
$ws.Range("K41").Value = "Mercy"

# --- Existing row 40 gained Recover/Death cells with a value of 0 ---
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0

# --- Update the view's selection to match the new active cell ---
$ws.Range("K38").Select()
